# Fix the "Extent" column (F) values: swap "WV, VA" to "VA, WV"
# for all rows on Sheet1 where that text appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(3,4,5,6,7,8,9,11,20,21,22,27,32,33,34,43,44,46,48,49,50,52,53,54,57,58,62,63,64,68,69,70,71,72,74,76)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = "VA, WV"  # Column F ("Extent")
}
